$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(4, 7).Value = "bp_percent"
$ws.Cells.Item(4, 8).Value = "cpx_percent"

$ws.Cells.Item(5, 7).Value = 100
$ws.Cells.Item(5, 8).Value = 100
$ws.Cells.Item(6, 7).Value = 46.333333333333336
$ws.Cells.Item(6, 8).Value = 104.71698113207547
$ws.Cells.Item(7, 7).Value = 23.333333333333332
$ws.Cells.Item(7, 8).Value = 141.50943396226415
$ws.Cells.Item(8, 7).Value = 100
$ws.Cells.Item(8, 8).Value = 100
$ws.Cells.Item(9, 7).Value = 6.219151036525172
$ws.Cells.Item(9, 8).Value = 115.58441558441558
$ws.Cells.Item(10, 7).Value = 5.923000987166831
$ws.Cells.Item(10, 8).Value = 100
$ws.Cells.Item(11, 7).Value = 100
$ws.Cells.Item(11, 8).Value = 100
$ws.Cells.Item(12, 7).Value = 63.116883116883116
$ws.Cells.Item(12, 8).Value = 157.33333333333334
$ws.Cells.Item(13, 7).Value = 12.987012987012987
$ws.Cells.Item(13, 8).Value = 90.66666666666667
$ws.Cells.Item(14, 7).Value = 100
$ws.Cells.Item(14, 8).Value = 100
$ws.Cells.Item(15, 7).Value = 24.673439767779392
$ws.Cells.Item(15, 8).Value = 531.8181818181818
$ws.Cells.Item(16, 7).Value = 7.547169811320755
$ws.Cells.Item(16, 8).Value = 504.5454545454545
$ws.Cells.Item(17, 7).Value = 100
$ws.Cells.Item(17, 8).Value = 100
$ws.Cells.Item(18, 7).Value = 21.035598705501616
$ws.Cells.Item(18, 8).Value = 74.80314960629921
$ws.Cells.Item(19, 7).Value = 7.44336569579288
$ws.Cells.Item(19, 8).Value = 107.08661417322836
$ws.Cells.Item(20, 7).Value = 100
$ws.Cells.Item(20, 8).Value = 100
$ws.Cells.Item(21, 7).Value = 21.073170731707318
$ws.Cells.Item(21, 8).Value = 46.666666666666664
$ws.Cells.Item(22, 7).Value = 18.634146341463413
$ws.Cells.Item(22, 8).Value = 105
$ws.Cells.Item(23, 7).Value = 100
$ws.Cells.Item(23, 8).Value = 100
$ws.Cells.Item(24, 7).Value = 17.4573055028463
$ws.Cells.Item(24, 8).Value = 88.54961832061069
$ws.Cells.Item(25, 7).Value = 12.903225806451612
$ws.Cells.Item(25, 8).Value = 78.62595419847328
$ws.Cells.Item(26, 7).Value = 100
$ws.Cells.Item(26, 8).Value = 100
$ws.Cells.Item(27, 7).Value = 99.86910994764398
$ws.Cells.Item(27, 8).Value = 68.35443037974683
$ws.Cells.Item(28, 7).Value = 44.24083769633508
$ws.Cells.Item(28, 8).Value = 82.27848101265822
$ws.Cells.Item(29, 7).Value = 100
$ws.Cells.Item(29, 8).Value = 100
$ws.Cells.Item(30, 7).Value = 53.43434343434343
$ws.Cells.Item(30, 8).Value = 89.23076923076923
$ws.Cells.Item(31, 7).Value = 10
$ws.Cells.Item(31, 8).Value = 107.6923076923077
$ws.Cells.Item(32, 7).Value = 100
$ws.Cells.Item(32, 8).Value = 100
$ws.Cells.Item(33, 7).Value = 23.674911660777386
$ws.Cells.Item(33, 8).Value = 115.38461538461539
$ws.Cells.Item(34, 7).Value = 21.554770318021202
$ws.Cells.Item(34, 8).Value = 105.76923076923076
$ws.Cells.Item(35, 7).Value = 100
$ws.Cells.Item(35, 8).Value = 100
$ws.Cells.Item(36, 7).Value = 30.192307692307693
$ws.Cells.Item(36, 8).Value = 108.47457627118644
$ws.Cells.Item(37, 7).Value = 14.807692307692308
$ws.Cells.Item(37, 8).Value = 100
$ws.Cells.Item(38, 7).Value = 100
$ws.Cells.Item(38, 8).Value = 100
$ws.Cells.Item(39, 7).Value = 20.88724584103512
$ws.Cells.Item(39, 8).Value = 86.92307692307692
$ws.Cells.Item(40, 7).Value = 12.199630314232902
$ws.Cells.Item(40, 8).Value = 97.6923076923077
$ws.Cells.Item(41, 7).Value = 100
$ws.Cells.Item(41, 8).Value = 100
$ws.Cells.Item(42, 7).Value = 23.343373493975903
$ws.Cells.Item(42, 8).Value = 172.64150943396228
$ws.Cells.Item(43, 7).Value = 12.349397590361447
$ws.Cells.Item(43, 8).Value = 204.7169811320755
$ws.Cells.Item(44, 7).Value = 100
$ws.Cells.Item(44, 8).Value = 100
$ws.Cells.Item(45, 7).Value = 54.67032967032967
$ws.Cells.Item(45, 8).Value = 93.24324324324324
$ws.Cells.Item(46, 7).Value = 33.79120879120879
$ws.Cells.Item(46, 8).Value = 82.43243243243244
$ws.Cells.Item(47, 7).Value = 100
$ws.Cells.Item(47, 8).Value = 100
$ws.Cells.Item(48, 7).Value = 9.877488514548238
$ws.Cells.Item(48, 8).Value = 91.22807017543859
$ws.Cells.Item(49, 7).Value = 5.359877488514548
$ws.Cells.Item(49, 8).Value = 133.33333333333334
$ws.Cells.Item(50, 7).Value = 100
$ws.Cells.Item(50, 8).Value = 100
$ws.Cells.Item(51, 7).Value = 24.84472049689441
$ws.Cells.Item(51, 8).Value = 150
$ws.Cells.Item(52, 7).Value = 15.320910973084887
$ws.Cells.Item(52, 8).Value = 113.5135135135135
$ws.Cells.Item(53, 7).Value = 100
$ws.Cells.Item(53, 8).Value = 100
$ws.Cells.Item(54, 7).Value = 9.491978609625669
$ws.Cells.Item(54, 8).Value = 81.35593220338983
$ws.Cells.Item(55, 7).Value = 8.02139037433155
$ws.Cells.Item(55, 8).Value = 117.79661016949152
$ws.Cells.Item(56, 7).Value = 100
$ws.Cells.Item(56, 8).Value = 100
$ws.Cells.Item(57, 7).Value = 23.716381418092908
$ws.Cells.Item(57, 8).Value = 84.53038674033148
$ws.Cells.Item(58, 7).Value = 19.070904645476773
$ws.Cells.Item(58, 8).Value = 78.45303867403314
$ws.Cells.Item(59, 7).Value = 100
$ws.Cells.Item(59, 8).Value = 100
$ws.Cells.Item(60, 7).Value = 51.785714285714285
$ws.Cells.Item(60, 8).Value = 95.40229885057474
$ws.Cells.Item(61, 7).Value = 9.226190476190476
$ws.Cells.Item(61, 8).Value = 189.6551724137931
$ws.Cells.Item(62, 7).Value = 100
$ws.Cells.Item(62, 8).Value = 100
$ws.Cells.Item(63, 7).Value = 21.409921671018278
$ws.Cells.Item(63, 8).Value = 42.3728813559322
$ws.Cells.Item(64, 7).Value = 14.099216710182768
$ws.Cells.Item(64, 8).Value = 105.9322033898305
$ws.Cells.Item(65, 7).Value = 100
$ws.Cells.Item(65, 8).Value = 100
$ws.Cells.Item(66, 7).Value = 15.954773869346734
$ws.Cells.Item(66, 8).Value = 156.66666666666666
$ws.Cells.Item(67, 7).Value = 11.4321608040201
$ws.Cells.Item(67, 8).Value = 83.33333333333333
$ws.Cells.Item(68, 7).Value = 100
$ws.Cells.Item(68, 8).Value = 100
$ws.Cells.Item(69, 7).Value = 7.387387387387387
$ws.Cells.Item(69, 8).Value = 127.08333333333334
$ws.Cells.Item(70, 7).Value = 7.207207207207207
$ws.Cells.Item(70, 8).Value = 0
$ws.Cells.Item(71, 7).Value = 100
$ws.Cells.Item(71, 8).Value = 100
$ws.Cells.Item(72, 7).Value = 31.645569620253166
$ws.Cells.Item(72, 8).Value = 132.09876543209876
$ws.Cells.Item(73, 7).Value = 23.417721518987342
$ws.Cells.Item(73, 8).Value = 100
$ws.Cells.Item(74, 7).Value = 100
$ws.Cells.Item(74, 8).Value = 100
$ws.Cells.Item(75, 7).Value = 67.37588652482269
$ws.Cells.Item(75, 8).Value = 74.66666666666667
$ws.Cells.Item(76, 7).Value = 10.99290780141844
$ws.Cells.Item(76, 8).Value = 81.33333333333333
$ws.Cells.Item(77, 7).Value = 100
$ws.Cells.Item(77, 8).Value = 100
$ws.Cells.Item(78, 7).Value = 28.29787234042553
$ws.Cells.Item(78, 8).Value = 123.4375
$ws.Cells.Item(79, 7).Value = 18.72340425531915
$ws.Cells.Item(79, 8).Value = 60.9375
$ws.Cells.Item(80, 7).Value = 100
$ws.Cells.Item(80, 8).Value = 100
$ws.Cells.Item(81, 7).Value = 25.358490566037737
$ws.Cells.Item(81, 8).Value = 121.875
$ws.Cells.Item(82, 7).Value = 14.11320754716981
$ws.Cells.Item(82, 8).Value = 128.12499999999997
$ws.Cells.Item(83, 7).Value = 100
$ws.Cells.Item(83, 8).Value = 100
$ws.Cells.Item(84, 7).Value = 9.956709956709958
$ws.Cells.Item(84, 8).Value = 85.71428571428571
$ws.Cells.Item(85, 7).Value = 7.575757575757576
$ws.Cells.Item(85, 8).Value = 101.19047619047619
$ws.Cells.Item(86, 7).Value = 100
$ws.Cells.Item(86, 8).Value = 100
$ws.Cells.Item(87, 7).Value = 19.167904903417533
$ws.Cells.Item(87, 8).Value = 111.50442477876105
$ws.Cells.Item(88, 7).Value = 10.846953937592868
$ws.Cells.Item(88, 8).Value = 112.38938053097344
$ws.Cells.Item(89, 7).Value = 100
$ws.Cells.Item(89, 8).Value = 100
$ws.Cells.Item(90, 7).Value = 70.33248081841433
$ws.Cells.Item(90, 8).Value = 98.7012987012987
$ws.Cells.Item(91, 7).Value = 34.271099744245525
$ws.Cells.Item(91, 8).Value = 103.8961038961039
$ws.Cells.Item(92, 7).Value = 100
$ws.Cells.Item(92, 8).Value = 100
$ws.Cells.Item(93, 7).Value = 26.34508348794063
$ws.Cells.Item(93, 8).Value = 117.64705882352942
$ws.Cells.Item(94, 7).Value = 12.059369202226344
$ws.Cells.Item(94, 8).Value = 75.49019607843138
$ws.Cells.Item(95, 7).Value = 100
$ws.Cells.Item(95, 8).Value = 0
$ws.Cells.Item(96, 7).Value = 21.566731141199227
$ws.Cells.Item(96, 8).Value = 0
$ws.Cells.Item(97, 7).Value = 11.218568665377177
$ws.Cells.Item(97, 8).Value = 0
$ws.Cells.Item(98, 7).Value = 100
$ws.Cells.Item(98, 8).Value = 100
$ws.Cells.Item(99, 7).Value = 53.099173553719005
$ws.Cells.Item(99, 8).Value = 117.54385964912281
$ws.Cells.Item(100, 7).Value = 11.363636363636363
$ws.Cells.Item(100, 8).Value = 112.28070175438596
$ws.Cells.Item(101, 7).Value = 100
$ws.Cells.Item(101, 8).Value = 100
$ws.Cells.Item(102, 7).Value = 63.16568047337278
$ws.Cells.Item(102, 8).Value = 120.45454545454544
$ws.Cells.Item(103, 7).Value = 23.668639053254438
$ws.Cells.Item(103, 8).Value = 143.18181818181816
$ws.Cells.Item(104, 7).Value = 100
$ws.Cells.Item(104, 8).Value = 100
$ws.Cells.Item(105, 7).Value = 81.26361655773421
$ws.Cells.Item(105, 8).Value = 70.2127659574468
$ws.Cells.Item(106, 7).Value = 35.947712418300654
$ws.Cells.Item(106, 8).Value = 225.531914893617
$ws.Cells.Item(107, 7).Value = 100
$ws.Cells.Item(107, 8).Value = 100
$ws.Cells.Item(108, 7).Value = 22.374429223744293
$ws.Cells.Item(108, 8).Value = 172.88135593220335
$ws.Cells.Item(109, 7).Value = 21.765601217656013
$ws.Cells.Item(109, 8).Value = 67.79661016949152
$ws.Cells.Item(110, 7).Value = 100
$ws.Cells.Item(110, 8).Value = 100
$ws.Cells.Item(111, 7).Value = 34.28165007112376
$ws.Cells.Item(111, 8).Value = 115.09433962264151
$ws.Cells.Item(112, 7).Value = 10.81081081081081
$ws.Cells.Item(112, 8).Value = 175.47169811320757
$ws.Cells.Item(113, 7).Value = 100
$ws.Cells.Item(113, 8).Value = 100
$ws.Cells.Item(114, 7).Value = 28.17824377457405
$ws.Cells.Item(114, 8).Value = 72.05882352941177
$ws.Cells.Item(115, 7).Value = 17.03800786369594
$ws.Cells.Item(115, 8).Value = 72.05882352941177
$ws.Cells.Item(116, 7).Value = 100
$ws.Cells.Item(116, 8).Value = 100
$ws.Cells.Item(117, 7).Value = 96.63608562691131
$ws.Cells.Item(117, 8).Value = 98.00000000000001
$ws.Cells.Item(118, 7).Value = 14.37308868501529
$ws.Cells.Item(118, 8).Value = 186.00000000000003
$ws.Cells.Item(119, 7).Value = 100
$ws.Cells.Item(119, 8).Value = 100
$ws.Cells.Item(120, 7).Value = 6.142857142857143
$ws.Cells.Item(120, 8).Value = 120.6896551724138
$ws.Cells.Item(121, 7).Value = 5.714285714285714
$ws.Cells.Item(121, 8).Value = 70.68965517241378
$ws.Cells.Item(122, 7).Value = 100
$ws.Cells.Item(122, 8).Value = 100
$ws.Cells.Item(123, 7).Value = 26.013513513513512
$ws.Cells.Item(123, 8).Value = 129.54545454545453
$ws.Cells.Item(124, 7).Value = 5.236486486486487
$ws.Cells.Item(124, 8).Value = 224.99999999999997
$ws.Cells.Item(125, 7).Value = 100
$ws.Cells.Item(125, 8).Value = 100
$ws.Cells.Item(126, 7).Value = 30.434782608695652
$ws.Cells.Item(126, 8).Value = 105.71428571428571
$ws.Cells.Item(127, 7).Value = 14.130434782608695
$ws.Cells.Item(127, 8).Value = 117.14285714285714
$ws.Cells.Item(128, 7).Value = 100
$ws.Cells.Item(128, 8).Value = 100
$ws.Cells.Item(129, 7).Value = 21.3953488372093
$ws.Cells.Item(129, 8).Value = 152.72727272727272
$ws.Cells.Item(130, 7).Value = 8.13953488372093
$ws.Cells.Item(130, 8).Value = 181.8181818181818
$ws.Cells.Item(131, 7).Value = 100
$ws.Cells.Item(131, 8).Value = 100
$ws.Cells.Item(132, 7).Value = 47.08029197080292
$ws.Cells.Item(132, 8).Value = 105.26315789473684
$ws.Cells.Item(133, 7).Value = 15.693430656934307
$ws.Cells.Item(133, 8).Value = 88.7218045112782
$ws.Cells.Item(134, 7).Value = 100
$ws.Cells.Item(134, 8).Value = 100
$ws.Cells.Item(135, 7).Value = 55.89041095890411
$ws.Cells.Item(135, 8).Value = 127.45098039215686
$ws.Cells.Item(136, 7).Value = 27.945205479452056
$ws.Cells.Item(136, 8).Value = 76.47058823529412
$ws.Cells.Item(137, 7).Value = 100
$ws.Cells.Item(137, 8).Value = 100
$ws.Cells.Item(138, 7).Value = 21.153846153846153
$ws.Cells.Item(138, 8).Value = 80.3921568627451
$ws.Cells.Item(139, 7).Value = 17.032967032967033
$ws.Cells.Item(139, 8).Value = 130.06535947712416
$ws.Cells.Item(140, 7).Value = 100
$ws.Cells.Item(140, 8).Value = 100
$ws.Cells.Item(141, 7).Value = 57.429718875502004
$ws.Cells.Item(141, 8).Value = 175
$ws.Cells.Item(142, 7).Value = 22.89156626506024
$ws.Cells.Item(142, 8).Value = 182.35294117647058
$ws.Cells.Item(143, 7).Value = 100
$ws.Cells.Item(143, 8).Value = 100
$ws.Cells.Item(144, 7).Value = 15.15748031496063
$ws.Cells.Item(144, 8).Value = 167.5925925925926
$ws.Cells.Item(145, 7).Value = 5.905511811023622
$ws.Cells.Item(145, 8).Value = 109.25925925925925
$ws.Cells.Item(146, 7).Value = 100
$ws.Cells.Item(146, 8).Value = 100
$ws.Cells.Item(147, 7).Value = 47.44897959183673
$ws.Cells.Item(147, 8).Value = 147.72727272727272
$ws.Cells.Item(148, 7).Value = 22.448979591836736
$ws.Cells.Item(148, 8).Value = 219.3181818181818
$ws.Cells.Item(149, 7).Value = 100
$ws.Cells.Item(149, 8).Value = 100
$ws.Cells.Item(150, 7).Value = 72.85714285714286
$ws.Cells.Item(150, 8).Value = 97.82608695652175
$ws.Cells.Item(151, 7).Value = 24.285714285714285
$ws.Cells.Item(151, 8).Value = 100
$ws.Cells.Item(152, 7).Value = 100
$ws.Cells.Item(152, 8).Value = 100
$ws.Cells.Item(153, 7).Value = 42.391304347826086
$ws.Cells.Item(153, 8).Value = 312.44635193133047
$ws.Cells.Item(154, 7).Value = 14.565217391304348
$ws.Cells.Item(154, 8).Value = 54.93562231759657
$ws.Cells.Item(155, 7).Value = 100
$ws.Cells.Item(155, 8).Value = 100
$ws.Cells.Item(156, 7).Value = 84.12698412698413
$ws.Cells.Item(156, 8).Value = 92.59259259259258
$ws.Cells.Item(157, 7).Value = 7.936507936507937
$ws.Cells.Item(157, 8).Value = 101.85185185185185
$ws.Cells.Item(158, 7).Value = 100
$ws.Cells.Item(158, 8).Value = 100
$ws.Cells.Item(159, 7).Value = 8.474576271186441
$ws.Cells.Item(159, 8).Value = 119.56521739130434
$ws.Cells.Item(160, 7).Value = 4.745762711864407
$ws.Cells.Item(160, 8).Value = 94.92753623188405
$ws.Cells.Item(161, 7).Value = 100
$ws.Cells.Item(161, 8).Value = 100
$ws.Cells.Item(162, 7).Value = 48.86363636363637
$ws.Cells.Item(162, 8).Value = 103.96039603960396
$ws.Cells.Item(163, 7).Value = 37.5
$ws.Cells.Item(163, 8).Value = 139.6039603960396
$ws.Cells.Item(164, 7).Value = 100
$ws.Cells.Item(164, 8).Value = 100
$ws.Cells.Item(165, 7).Value = 29.661016949152543
$ws.Cells.Item(165, 8).Value = 102.59740259740259
$ws.Cells.Item(166, 7).Value = 17.51412429378531
$ws.Cells.Item(166, 8).Value = 88.31168831168831
$ws.Cells.Item(167, 7).Value = 100
$ws.Cells.Item(167, 8).Value = 100
$ws.Cells.Item(168, 7).Value = 57.38636363636363
$ws.Cells.Item(168, 8).Value = 141.1764705882353
$ws.Cells.Item(169, 7).Value = 14.772727272727273
$ws.Cells.Item(169, 8).Value = 138.23529411764707
$ws.Cells.Item(170, 7).Value = 100
$ws.Cells.Item(170, 8).Value = 100
$ws.Cells.Item(171, 7).Value = 61.904761904761905
$ws.Cells.Item(171, 8).Value = 84.13793103448276
$ws.Cells.Item(172, 7).Value = 19.841269841269842
$ws.Cells.Item(172, 8).Value = 64.48275862068965
$ws.Cells.Item(173, 7).Value = 100
$ws.Cells.Item(173, 8).Value = 100
$ws.Cells.Item(174, 7).Value = 31.16883116883117
$ws.Cells.Item(174, 8).Value = 154.71698113207546
$ws.Cells.Item(175, 7).Value = 18.181818181818183
$ws.Cells.Item(175, 8).Value = 103.77358490566039
$ws.Cells.Item(176, 7).Value = 100
$ws.Cells.Item(176, 8).Value = 100
$ws.Cells.Item(177, 7).Value = 33.333333333333336
$ws.Cells.Item(177, 8).Value = 148.66666666666666
$ws.Cells.Item(178, 7).Value = 33.333333333333336
$ws.Cells.Item(178, 8).Value = 94
$ws.Cells.Item(179, 7).Value = 100
$ws.Cells.Item(179, 8).Value = 100
$ws.Cells.Item(180, 7).Value = 68.35443037974683
$ws.Cells.Item(180, 8).Value = 217.77777777777777
$ws.Cells.Item(181, 7).Value = 48.10126582278481
$ws.Cells.Item(181, 8).Value = 233.33333333333334
$ws.Cells.Item(182, 7).Value = 100
$ws.Cells.Item(182, 8).Value = 100
$ws.Cells.Item(183, 7).Value = 31.914893617021278
$ws.Cells.Item(183, 8).Value = 75
$ws.Cells.Item(184, 7).Value = 17.02127659574468
$ws.Cells.Item(184, 8).Value = 91.34615384615384
$ws.Cells.Item(185, 7).Value = 100
$ws.Cells.Item(185, 8).Value = 100
$ws.Cells.Item(186, 7).Value = 75.625
$ws.Cells.Item(186, 8).Value = 66.42335766423358
$ws.Cells.Item(187, 7).Value = 18.75
$ws.Cells.Item(187, 8).Value = 129.56204379562044
$ws.Cells.Item(188, 7).Value = 100
$ws.Cells.Item(188, 8).Value = 100
$ws.Cells.Item(189, 7).Value = 38.9937106918239
$ws.Cells.Item(189, 8).Value = 177.7251184834123
$ws.Cells.Item(190, 7).Value = 30.81761006289308
$ws.Cells.Item(190, 8).Value = 86.7298578199052
$ws.Cells.Item(191, 7).Value = 100
$ws.Cells.Item(191, 8).Value = 100
$ws.Cells.Item(192, 7).Value = 91.83673469387755
$ws.Cells.Item(192, 8).Value = 217.92114695340504
$ws.Cells.Item(193, 7).Value = 55.10204081632653
$ws.Cells.Item(193, 8).Value = 143.01075268817206
$ws.Cells.Item(194, 7).Value = 100
$ws.Cells.Item(194, 8).Value = 100
$ws.Cells.Item(195, 7).Value = 80
$ws.Cells.Item(195, 8).Value = 33.557046979865774
$ws.Cells.Item(196, 7).Value = 20
$ws.Cells.Item(196, 8).Value = 30.872483221476507
$ws.Cells.Item(197, 7).Value = 100
$ws.Cells.Item(197, 8).Value = 100
$ws.Cells.Item(198, 7).Value = 50.40650406504065
$ws.Cells.Item(198, 8).Value = 65.21739130434783
$ws.Cells.Item(199, 7).Value = 13.821138211382113
$ws.Cells.Item(199, 8).Value = 0
$ws.Cells.Item(200, 7).Value = 100
$ws.Cells.Item(200, 8).Value = 100
$ws.Cells.Item(201, 7).Value = 62.5
$ws.Cells.Item(201, 8).Value = 82.95819935691318
$ws.Cells.Item(202, 7).Value = 57.5
$ws.Cells.Item(202, 8).Value = 90.67524115755627
$ws.Cells.Item(203, 7).Value = 100
$ws.Cells.Item(203, 8).Value = 100
$ws.Cells.Item(204, 7).Value = 85.4368932038835
$ws.Cells.Item(204, 8).Value = 173.41772151898732
$ws.Cells.Item(205, 7).Value = 27.184466019417474
$ws.Cells.Item(205, 8).Value = 50.63291139240506
$ws.Cells.Item(206, 7).Value = 100
$ws.Cells.Item(206, 8).Value = 100
$ws.Cells.Item(207, 7).Value = 81.63265306122449
$ws.Cells.Item(207, 8).Value = 93.33333333333333
$ws.Cells.Item(208, 7).Value = 32.6530612244898
$ws.Cells.Item(208, 8).Value = 107.45098039215686
$ws.Cells.Item(209, 7).Value = 100
$ws.Cells.Item(209, 8).Value = 100
$ws.Cells.Item(210, 7).Value = 96.875
$ws.Cells.Item(210, 8).Value = 124.4019138755981
$ws.Cells.Item(211, 7).Value = 40.625
$ws.Cells.Item(211, 8).Value = 107.6555023923445
$ws.Cells.Item(212, 7).Value = 100
$ws.Cells.Item(212, 8).Value = 100
$ws.Cells.Item(213, 7).Value = 100
$ws.Cells.Item(213, 8).Value = 76.37795275590551
$ws.Cells.Item(214, 7).Value = 34.04255319148936
$ws.Cells.Item(214, 8).Value = 70.07874015748033
$ws.Cells.Item(215, 7).Value = 100
$ws.Cells.Item(215, 8).Value = 100
$ws.Cells.Item(216, 7).Value = 18.27956989247312
$ws.Cells.Item(216, 8).Value = 125.75757575757578
$ws.Cells.Item(217, 7).Value = 17.204301075268816
$ws.Cells.Item(217, 8).Value = 0
$ws.Cells.Item(218, 7).Value = 100
$ws.Cells.Item(218, 8).Value = 100
$ws.Cells.Item(219, 7).Value = 43.13725490196079
$ws.Cells.Item(219, 8).Value = 119.9124726477024
$ws.Cells.Item(220, 7).Value = 9.803921568627452
$ws.Cells.Item(220, 8).Value = 180.08752735229757
$ws.Cells.Item(221, 7).Value = 100
$ws.Cells.Item(221, 8).Value = 100
$ws.Cells.Item(222, 7).Value = 47.82608695652174
$ws.Cells.Item(222, 8).Value = 75.47169811320755
$ws.Cells.Item(223, 7).Value = 15.942028985507246
$ws.Cells.Item(223, 8).Value = 89.937106918239
$ws.Cells.Item(224, 7).Value = 100
$ws.Cells.Item(224, 8).Value = 100
$ws.Cells.Item(225, 7).Value = 42.10526315789474
$ws.Cells.Item(225, 8).Value = 152.6530612244898
$ws.Cells.Item(226, 7).Value = 26.31578947368421
$ws.Cells.Item(226, 8).Value = 148.16326530612244
$ws.Cells.Item(227, 7).Value = 100
$ws.Cells.Item(227, 8).Value = 0
$ws.Cells.Item(228, 7).Value = 46.42857142857143
$ws.Cells.Item(228, 8).Value = 0
$ws.Cells.Item(229, 7).Value = 12.5
$ws.Cells.Item(229, 8).Value = 0
$ws.Cells.Item(230, 7).Value = 100
$ws.Cells.Item(230, 8).Value = 100
$ws.Cells.Item(231, 7).Value = 37.77777777777778
$ws.Cells.Item(231, 8).Value = 183.1831831831832
$ws.Cells.Item(232, 7).Value = 11.11111111111111
$ws.Cells.Item(232, 8).Value = 146.54654654654655
$ws.Cells.Item(233, 7).Value = 100
$ws.Cells.Item(233, 8).Value = 100
$ws.Cells.Item(234, 7).Value = 26.31578947368421
$ws.Cells.Item(234, 8).Value = 65.9090909090909
$ws.Cells.Item(235, 7).Value = 18.42105263157895
$ws.Cells.Item(235, 8).Value = 75.75757575757576
$ws.Cells.Item(236, 7).Value = 100
$ws.Cells.Item(236, 8).Value = 100
$ws.Cells.Item(237, 7).Value = 75
$ws.Cells.Item(237, 8).Value = 87.56756756756756
$ws.Cells.Item(238, 7).Value = 62.5
$ws.Cells.Item(238, 8).Value = 98.64864864864865
$ws.Cells.Item(239, 7).Value = 100
$ws.Cells.Item(239, 8).Value = 100
$ws.Cells.Item(240, 7).Value = 45.90163934426229
$ws.Cells.Item(240, 8).Value = 108.33333333333333
$ws.Cells.Item(241, 7).Value = 8.19672131147541
$ws.Cells.Item(241, 8).Value = 122.5
$ws.Cells.Item(242, 7).Value = 100
$ws.Cells.Item(242, 8).Value = 100
$ws.Cells.Item(243, 7).Value = 34.61538461538461
$ws.Cells.Item(243, 8).Value = 136.7983367983368
$ws.Cells.Item(244, 7).Value = 21.153846153846153
$ws.Cells.Item(244, 8).Value = 79.62577962577961
